$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers) ---
# Copy the existing header style (bold, border, centered) from B1 onto the
# new header cells F1:K1 before writing their text.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("F1:K1").PasteSpecial(-4122) | Out-Null

$ws.Range("B1").Value = "L-1NANE"
$ws.Range("C1").Value = "L3NAYE"
$ws.Range("D1").Value = "L0NAYE"
$ws.Range("E1").Value = "L3NANE"
$ws.Range("F1").Value = "L-1YANE"
$ws.Range("G1").Value = "L1NAYE"
$ws.Range("H1").Value = "L2NANE"
$ws.Range("I1").Value = "L2NAYE"
$ws.Range("J1").Value = "L3YANE"
$ws.Range("K1").Value = "L2YANE"

# --- Row 2 ("Steps On Failure") ---
$ws.Range("B2").Value = 718.8
$ws.Range("C2").Value = 234.72
$ws.Range("D2").Value = 72.61
$ws.Range("E2").Value = 379.01
$ws.Range("F2").Value = 870.4
$ws.Range("G2").Value = 64.81999999999999

# H2 / K2 hold the text "0" (not the number 0). Force text by formatting the
# cell as Text, writing the value, then clearing the format back off again
# so no lingering style index remains on the cell.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "0"
$ws.Range("H2").ClearFormats()

$ws.Range("I2").Value = 24
$ws.Range("J2").Value = 569.3

$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "0"
$ws.Range("K2").ClearFormats()

# --- Row 3 ("Steps On Success") ---
$ws.Range("B3").Value = 413.66
$ws.Range("C3").Value = 281.85
$ws.Range("D3").Value = 81.98999999999999
$ws.Range("E3").Value = 276.15
$ws.Range("F3").Value = 445.07
$ws.Range("G3").Value = 60.08
$ws.Range("H3").Value = 126.025
$ws.Range("I3").Value = 146.18
$ws.Range("J3").Value = 277.6
$ws.Range("K3").Value = 127.7549019607843

Write-Output "edit applied"
